$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2364341085271318
$ws.Range("C2").Value = 0.4883720930232558
$ws.Range("J2").Value = 0.01937984496124031
$ws.Range("P2").Value = 0.1705426356589147
$ws.Range("S2").Value = 0.08527131782945736
$ws.Range("B3").Value = 0.01574803149606299
$ws.Range("C3").Value = 0.01574803149606299
$ws.Range("J3").Value = 0.03149606299212598
$ws.Range("P3").Value = 0.7322834645669292
$ws.Range("S3").Value = 0.2047244094488189
$ws.Range("J4").Value = 0.02857142857142857
$ws.Range("P4").Value = 0.6285714285714286
$ws.Range("S4").Value = 0.3428571428571429
$ws.Range("B6").Value = 0.02339181286549707
$ws.Range("D6").Value = 0.01754385964912281
$ws.Range("F6").Value = 0.04678362573099415
$ws.Range("J6").Value = 0.239766081871345
$ws.Range("O6").Value = 0.02339181286549707
$ws.Range("Q6").Value = 0.152046783625731
$ws.Range("R6").Value = 0.09941520467836257
$ws.Range("S6").Value = 0.3976608187134503
$ws.Range("B7").Value = 0.08823529411764706
$ws.Range("D7").Value = 0.01764705882352941
$ws.Range("F7").Value = 0.03529411764705882
$ws.Range("J7").Value = 0.1588235294117647
$ws.Range("O7").Value = 0.005882352941176471
$ws.Range("Q7").Value = 0.2176470588235294
$ws.Range("R7").Value = 0.07647058823529412
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.08635097493036212
$ws.Range("D8").Value = 0.01392757660167131
$ws.Range("F8").Value = 0.06128133704735376
$ws.Range("J8").Value = 0.1114206128133705
$ws.Range("O8").Value = 0.01949860724233983
$ws.Range("Q8").Value = 0.1643454038997214
$ws.Range("R8").Value = 0.1002785515320334
$ws.Range("S8").Value = 0.4428969359331476
$ws.Range("B9").Value = 0.1040462427745665
$ws.Range("D9").Value = 0.01734104046242774
$ws.Range("E9").Value = 0.005780346820809248
$ws.Range("F9").Value = 0.04046242774566474
$ws.Range("J9").Value = 0.06358381502890173
$ws.Range("O9").Value = 0.0115606936416185
$ws.Range("Q9").Value = 0.1734104046242775
$ws.Range("R9").Value = 0.1040462427745665
$ws.Range("S9").Value = 0.4797687861271676
$ws.Range("B10").Value = 0.120952380952381
$ws.Range("D10").Value = 0.02190476190476191
$ws.Range("E10").Value = 0.001904761904761905
$ws.Range("F10").Value = 0.07809523809523809
$ws.Range("J10").Value = 0.1066666666666667
$ws.Range("O10").Value = 0.01333333333333333
$ws.Range("Q10").Value = 0.1819047619047619
$ws.Range("R10").Value = 0.08476190476190476
$ws.Range("S10").Value = 0.3904761904761905
$ws.Range("G11").Value = 0.1259541984732824
$ws.Range("J11").Value = 0.0916030534351145
$ws.Range("K11").Value = 0.1755725190839695
$ws.Range("L11").Value = 0.583969465648855
$ws.Range("S11").Value = 0.02290076335877863
$ws.Range("G12").Value = 0.7453416149068323
$ws.Range("J12").Value = 0.1987577639751553
$ws.Range("L12").Value = 0.02484472049689441
$ws.Range("S12").Value = 0.03105590062111801
$ws.Range("G13").Value = 0.6216216216216216
$ws.Range("J13").Value = 0.2972972972972973
$ws.Range("S13").Value = 0.08108108108108109
$ws.Range("F15").Value = 0.01485148514851485
$ws.Range("H15").Value = 0.1683168316831683
$ws.Range("I15").Value = 0.08415841584158416
$ws.Range("J15").Value = 0.3415841584158416
$ws.Range("K15").Value = 0.08415841584158416
$ws.Range("M15").Value = 0.004950495049504951
$ws.Range("O15").Value = 0.06930693069306931
$ws.Range("S15").Value = 0.2326732673267327
$ws.Range("F16").Value = 0.006329113924050633
$ws.Range("H16").Value = 0.1518987341772152
$ws.Range("I16").Value = 0.1012658227848101
$ws.Range("J16").Value = 0.379746835443038
$ws.Range("K16").Value = 0.1582278481012658
$ws.Range("M16").Value = 0.0189873417721519
$ws.Range("O16").Value = 0.06329113924050633
$ws.Range("S16").Value = 0.120253164556962
$ws.Range("F17").Value = 0.02058823529411765
$ws.Range("H17").Value = 0.1823529411764706
$ws.Range("I17").Value = 0.1117647058823529
$ws.Range("J17").Value = 0.3882352941176471
$ws.Range("K17").Value = 0.09411764705882353
$ws.Range("M17").Value = 0.01764705882352941
$ws.Range("O17").Value = 0.05882352941176471
$ws.Range("S17").Value = 0.1264705882352941
$ws.Range("F18").Value = 0.01169590643274854
$ws.Range("H18").Value = 0.2046783625730994
$ws.Range("I18").Value = 0.06432748538011696
$ws.Range("J18").Value = 0.4210526315789473
$ws.Range("K18").Value = 0.09941520467836257
$ws.Range("M18").Value = 0.01169590643274854
$ws.Range("O18").Value = 0.07017543859649122
$ws.Range("S18").Value = 0.1169590643274854
$ws.Range("F19").Value = 0.01376146788990826
$ws.Range("H19").Value = 0.1871559633027523
$ws.Range("I19").Value = 0.0834862385321101
$ws.Range("J19").Value = 0.3880733944954128
$ws.Range("K19").Value = 0.1119266055045872
$ws.Range("M19").Value = 0.02385321100917431
$ws.Range("O19").Value = 0.08623853211009175
$ws.Range("S19").Value = 0.1055045871559633
